# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Choclo" (Choclero, Primera) above
# the existing row 99, shifting the rest of the data block down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 99; everything currently on/after row 99 (down to
# the former last row 119) shifts down to rows 100..120.
$ws.Rows("99:99").Insert()

# Populate the newly-inserted row 99 with the new weekly record.
$ws.Cells.Item(99, 1).Value = 2
$ws.Cells.Item(99, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(99, 3).Value = "Coquimbo"
$ws.Cells.Item(99, 4).Value = 44644
$ws.Cells.Item(99, 5).Value = 4
$ws.Cells.Item(99, 6).Value = 100112024
$ws.Cells.Item(99, 7).Value = "Choclo"
$ws.Cells.Item(99, 8).Value = "Choclero"
$ws.Cells.Item(99, 9).Value = "Primera"
$ws.Cells.Item(99, 10).Value = 40000
$ws.Cells.Item(99, 11).Value = 200
$ws.Cells.Item(99, 12).Value = 230
$ws.Cells.Item(99, 13).Value = 215
$ws.Cells.Item(99, 14).Value = "`$/unidad"
$ws.Cells.Item(99, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(99, 16).Value = 215
$ws.Cells.Item(99, 17).Value = 1
$ws.Cells.Item(99, 18).Value = "Hortaliza"
